$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E3").Value = 36
$ws.Range("F3").Value = 28
$ws.Range("H3").Value = 28
$ws.Range("F4").Value = 27
$ws.Range("H4").Value = 27
$ws.Range("E6").Value = 43
$ws.Range("E10").Value = 474
$ws.Range("F10").Value = 237
$ws.Range("H10").Value = 237
$ws.Range("E11").Value = 318
$ws.Range("F11").Value = 179
$ws.Range("H11").Value = 179
$ws.Range("E12").Value = 464
$ws.Range("E13").Value = 119
$ws.Range("E15").Value = 153
$ws.Range("E16").Value = 189
$ws.Range("F17").Value = 46
$ws.Range("H17").Value = 46
$ws.Range("E18").Value = 51
$ws.Range("F18").Value = 26
$ws.Range("H18").Value = 26
$ws.Range("E22").Value = 159
$ws.Range("F22").Value = 85
$ws.Range("H22").Value = 85
$ws.Range("F23").Value = 85
$ws.Range("H23").Value = 85
$ws.Range("E24").Value = 199
$ws.Range("F24").Value = 108
$ws.Range("H24").Value = 108
$ws.Range("F25").Value = 120
$ws.Range("H25").Value = 120
$ws.Range("E27").Value = 309
$ws.Range("F27").Value = 151
$ws.Range("H27").Value = 151
$ws.Range("E28").Value = 189
$ws.Range("F28").Value = 71
$ws.Range("H28").Value = 71
$ws.Range("E29").Value = 161
$ws.Range("F29").Value = 91
$ws.Range("H29").Value = 91
$ws.Range("F30").Value = 117
$ws.Range("H30").Value = 117
$ws.Range("E32").Value = 175
$ws.Range("E33").Value = 272
$ws.Range("F33").Value = 139
$ws.Range("H33").Value = 139
$ws.Range("E34").Value = 203
$ws.Range("F34").Value = 132
$ws.Range("H34").Value = 132
$ws.Range("E35").Value = 140
$ws.Range("F35").Value = 87
$ws.Range("H35").Value = 87
$ws.Range("E38").Value = 86
$ws.Range("E39").Value = 173
$ws.Range("F39").Value = 82
$ws.Range("H39").Value = 82
$ws.Range("E40").Value = 244
$ws.Range("F40").Value = 113
$ws.Range("H40").Value = 113
$ws.Range("E41").Value = 378
$ws.Range("F41").Value = 174
$ws.Range("H41").Value = 174
$ws.Range("E42").Value = 349
$ws.Range("F42").Value = 189
$ws.Range("H42").Value = 189
$ws.Range("E44").Value = 296
$ws.Range("E45").Value = 133
$ws.Range("F45").Value = 68
$ws.Range("H45").Value = 68
$ws.Range("E46").Value = 296
$ws.Range("F46").Value = 164
$ws.Range("H46").Value = 164
$ws.Range("E47").Value = 422
$ws.Range("F47").Value = 208
$ws.Range("H47").Value = 208
$ws.Range("E48").Value = 192
$ws.Range("F48").Value = 83
$ws.Range("H48").Value = 83
$ws.Range("E49").Value = 274
$ws.Range("F49").Value = 117
$ws.Range("H49").Value = 117
$ws.Range("E50").Value = 232
$ws.Range("F50").Value = 107
$ws.Range("H50").Value = 107
$ws.Range("E51").Value = 221
$ws.Range("F51").Value = 93
$ws.Range("H51").Value = 93
